$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell $ws "D2" '43.321.79'
Set-TextCell $ws "E2" '  -1.88%  '
Set-TextCell $ws "D3" '2.236.65'
Set-TextCell $ws "E3" '  -1.84%  '
Set-TextCell $ws "E4" '  +0.12%  '
Set-TextCell $ws "D5" '230.25'
Set-TextCell $ws "E5" '  -1.16%  '
Set-TextCell $ws "D6" '0.639'
Set-TextCell $ws "E6" '  -0.78%  '
Set-TextCell $ws "D7" '63.11'
Set-TextCell $ws "E7" '  -4.10%  '
Set-TextCell $ws "E8" '  +0.03%  '
Set-TextCell $ws "D9" '0.441'
Set-TextCell $ws "E9" '  +0.96%  '
Set-TextCell $ws "D10" '0.0952'
Set-TextCell $ws "E10" '  -7.19%  '
Set-TextCell $ws "D11" '56.63'
Set-TextCell $ws "E11" '  -1.61%  '
Set-TextCell $ws "D12" '27.43'
Set-TextCell $ws "E12" '  +4.55%  '
Set-TextCell $ws "E13" '  -0.27%  '
Set-TextCell $ws "D14" '2.573.79'
Set-TextCell $ws "E14" '  -1.62%  '
Set-TextCell $ws "D15" '15.39'
Set-TextCell $ws "E15" '  -3.73%  '
Set-TextCell $ws "D16" '6.04'
Set-TextCell $ws "E16" '  +0.31%  '
Set-TextCell $ws "D17" '0.824'
Set-TextCell $ws "E17" '  -1.66%  '
Set-TextCell $ws "D18" '2.242.80'
Set-TextCell $ws "E18" '  -1.46%  '
Set-TextCell $ws "D19" '43.210.99'
Set-TextCell $ws "E19" '  -1.65%  '
Set-TextCell $ws "D20" [string]::Concat('0.0', [string][char]8323, '0961')
Set-TextCell $ws "E20" '  -3.02%  '
Set-TextCell $ws "D21" '72.72'
Set-TextCell $ws "E21" '  -1.84%  '
Set-TextCell $ws "D22" '6.06'
Set-TextCell $ws "E22" '  -1.20%  '
Set-TextCell $ws "D23" '245.81'
Set-TextCell $ws "E23" '  -5.42%  '
Set-TextCell $ws "D24" '0.999'
Set-TextCell $ws "E24" '  -0.01%  '
Set-TextCell $ws "D25" '3.66'
Set-TextCell $ws "E25" '  +30.28%  '
Set-TextCell $ws "E26" '  -3.17%  '
Set-TextCell $ws "D27" '2.27'
Set-TextCell $ws "E27" '  -1.90%  '
Set-TextCell $ws "D28" '9.72'
Set-TextCell $ws "E28" '  -5.02%  '
Set-TextCell $ws "D29" '172.80'
Set-TextCell $ws "E29" '  +0.54%  '
Set-TextCell $ws "D30" '21.41'
Set-TextCell $ws "E30" '  +1.73%  '
Set-TextCell $ws "D31" '0.129'
Set-TextCell $ws "E31" '  -6.82%  '
Set-TextCell $ws "D32" '1.41'
Set-TextCell $ws "E32" '  -1.73%  '
Set-TextCell $ws "D33" '0.124'
Set-TextCell $ws "E33" '  -0.03%  '
Set-TextCell $ws "D34" '4.86'
Set-TextCell $ws "E34" '  +2.38%  '
Set-TextCell $ws "D35" '0.0670'
Set-TextCell $ws "E35" '  -3.31%  '
Set-TextCell $ws "D36" '4.85'
Set-TextCell $ws "E36" '  -3.75%  '
Set-TextCell $ws "D37" '3.57'
Set-TextCell $ws "E37" '  -7.88%  '
Set-TextCell $ws "D38" '6.27'
Set-TextCell $ws "E38" '  -8.38%  '
Set-TextCell $ws "E39" '  -4.34%  '
Set-TextCell $ws "D40" '0.0249'
Set-TextCell $ws "E40" '  -0.11%  '
Set-TextCell $ws "E41" '  -0.05%  '
Set-TextCell $ws "D42" '8.59'
Set-TextCell $ws "E42" '  +2.69%  '
Set-TextCell $ws "D43" '4.46'
Set-TextCell $ws "E43" '  -0.69%  '
Set-TextCell $ws "D44" '16.95'
Set-TextCell $ws "E44" '  -3.79%  '
Set-TextCell $ws "D45" '96.14'
Set-TextCell $ws "E45" '  -2.05%  '
Set-TextCell $ws "B46" 'TerraClassic'
Set-TextCell $ws "C46" 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
Set-TextCell $ws "D46" '0.000209'
Set-TextCell $ws "E46" '  +2.91%  '
Set-TextCell $ws "B47" 'Cronos'
Set-TextCell $ws "C47" 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws "D47" '0.0937'
Set-TextCell $ws "E47" '  -4.11%  '
Set-TextCell $ws "E48" '  -2.28%  '
Set-TextCell $ws "D49" '1.448.37'
Set-TextCell $ws "E49" '  -1.66%  '
Set-TextCell $ws "B50" 'NEARProtocol'
Set-TextCell $ws "C50" 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws "D50" '2.28'
Set-TextCell $ws "E50" '  -3.95%  '
Set-TextCell $ws "B51" 'Celestia'
Set-TextCell $ws "C51" 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextCell $ws "D51" '9.80'
Set-TextCell $ws "E51" '  -1.29%  '

"done"